$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.196.70'
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("D3").Value = '3.319.32'
$ws.Range("E3").Value = '  +0.05%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '582.84'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +3.21%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '182.93'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.69%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +2.75%  '
$ws.Range("D9").Value = '3.312.75'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("E11").Value = '  +1.16%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '46.34'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("E13").Value = '  +3.78%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '636.35'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +6.81%  '
$ws.Range("D15").Value = '3.853.06'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '68.269.09'
$ws.Range("E17").Value = '  +3.25%  '
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").Value = '3.322.45'
$ws.Range("E19").Value = '  +0.15%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.69'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.16%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '10.93'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.37%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.902'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '17.64'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.41%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.09'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.78%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '97.02'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.88%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.99'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.24%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.77'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.91%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.58'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.76%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '32.42'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +5.37%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '8.59'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("E31").Value = '  +0.94%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '604.96'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +8.14%  '
$ws.Range("D33").Value = '3.947.85'
$ws.Range("E33").Value = '  +3.78%  '
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("E35").Value = '  +1.82%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.51'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -4.01%  '
$ws.Range("E37").Value = '  +0.04%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '55.89'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("E39").Value = '  +3.93%  '
$ws.Range("E40").Value = '  +1.13%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.68'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +3.09%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '32.63'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.85%  '
$ws.Range("E43").Value = '  +0.09%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.39'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  +1.63%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0415'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("E47").Value = '  +14.86%  '
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("E49").Value = '  +0.58%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.55'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.16%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '131.03'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.21%  '
